# Add file upload functionality
# Appends 3 new log rows (124-126) to each of the 4 worksheets, mirroring
# the existing row layout (time, length, id, actual length, checksum and
# their decimal counterparts).

$wb = $excel.ActiveWorkbook

# Date/time values (Excel serial date numbers) shared by the new rows
# across all sheets.
$dt124 = [double]"45910.49049768518"
$dt125 = [double]"45911.49237268518"
$dt126 = [double]"45912.49243055555"

function Add-LogRow {
    param(
        $ws,
        [int]$row,
        [double]$dateVal,
        [string]$b,
        [string]$c,
        [string]$d,
        [string]$e,
        [double]$f,
        [double]$g,
        [double]$h,
        [double]$i
    )

    $ws.Range("A$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("A$row").Value = $dateVal
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
    $ws.Range("I$row").Value = $i
}

# ---- Sheet 1: FE_LFT_#1 ----
$ws1 = $wb.Worksheets.Item(1)
$g1 = [double]"7.598631275147109e+23"
Add-LogRow $ws1 124 $dt124 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xE0" "0xf" 380 $g1 228 15
Add-LogRow $ws1 125 $dt125 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xE0" "0xf" 380 $g1 228 15
Add-LogRow $ws1 126 $dt126 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xE0" "0xf" 380 $g1 228 15

# ---- Sheet 2: FE_LFT_#2 ----
$ws2 = $wb.Worksheets.Item(2)
$g2 = [double]"5.68432987514711e+23"
Add-LogRow $ws2 124 $dt124 "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xEC" "0xe" 400 $g2 240 14
Add-LogRow $ws2 125 $dt125 "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xEC" "0xe" 400 $g2 236 14
Add-LogRow $ws2 126 $dt126 "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xE8" "0xe" 400 $g2 236 14

# ---- Sheet 3: FE_PLT_#1 ----
$ws3 = $wb.Worksheets.Item(3)
$g3 = [double]"5.68631262647114e+23"
Add-LogRow $ws3 124 $dt124 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x57" "0x3" 110 $g3 87 3
Add-LogRow $ws3 125 $dt125 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x57" "0x3" 110 $g3 87 3
Add-LogRow $ws3 126 $dt126 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x57" "0x3" 110 $g3 87 3

# ---- Sheet 4: FE_PLT_#2 ----
$ws4 = $wb.Worksheets.Item(4)
$g4 = [double]"9.85046333984776e+23"
Add-LogRow $ws4 124 $dt124 "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x54" "0x3" 110 $g4 84 3
Add-LogRow $ws4 125 $dt125 "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x54" "0x3" 110 $g4 84 3
Add-LogRow $ws4 126 $dt126 "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x54" "0x3" 110 $g4 84 3

Write-Host "Added rows 124-126 to all 4 sheets"
